$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 29 and 30: coin identities swap (PEPE <-> Aptos)
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'6.05"
$ws.Range("E29").Value = "  +1.48%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0₃0700"
$ws.Range("E30").Value = "  +1.81%  "

# Price / Volume(1h) updates
$ws.Range("D2").Value = "'54.632.61"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "'2.278.32"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'504.18"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'128.55"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "'2.294.30"
$ws.Range("D10").Value = "'0.0964"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("D13").Value = "'4.91"
$ws.Range("E13").Value = "  +4.37%  "
$ws.Range("D14").Value = "'23.29"
$ws.Range("E14").Value = "  +3.41%  "
$ws.Range("D15").Value = "'2.683.07"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "'54.674.85"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "'2.278.31"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "'4.13"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "'306.61"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "'6.43"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "'60.40"
$ws.Range("E24").Value = "  -2.70%  "
$ws.Range("D25").Value = "'0.994"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("D27").Value = "'7.45"
$ws.Range("E27").Value = "  +1.74%  "
$ws.Range("D28").Value = "'170.72"
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("E32").Value = "  +3.16%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'17.93"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("D36").Value = "'0.911"
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").Value = "'3.78"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").Value = "'36.50"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("D40").Value = "'0.374"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D42").Value = "'3.38"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").Value = "'126.54"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("D45").Value = "'250.38"
$ws.Range("E45").Value = "  +4.43%  "
$ws.Range("D46").Value = "'0.0495"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "'0.373"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("E51").Value = "  +0.42%  "
